$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.415.64'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.893.56'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.81%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.55'
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("E6").Value = '  -4.27%  '
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.94'
$ws.Range("E8").Value = '  +8.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.352'
$ws.Range("E9").Value = '  -5.50%  '
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '13.11'
$ws.Range("E12").Value = '  +2.11%  '
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.903.44'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.412.47'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '73.79'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '247.52'
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.84'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.96'
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("E24").Value = '  +4.95%  '
$ws.Range("E25").Value = '  -10.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.46'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.46'
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.36'
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("E29").Value = '  -3.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.128.48'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  +8.14%  '
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.22'
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.857'
$ws.Range("E36").Value = '  -6.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.02'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("E38").Value = '  -20.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0682'
$ws.Range("E39").Value = '  +4.32%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.14'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.53'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.295.85'
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.37'
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("E46").Value = '  +6.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.74'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.19'
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.36'
$ws.Range("E50").Value = '  -5.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.28'
$ws.Range("E51").Value = '  -3.63%  '
